$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-13 Thursday" "2025-03-14 Friday"

Replace-Text "256×5=" "346×8="
Replace-Text "849×5=" "225×8="
Replace-Text "703×6=" "874×3="
Replace-Text "917×8=" "612×9="
Replace-Text "423×3=" "640×9="

Replace-Text "243×7=" "564×6="
Replace-Text "401×3=" "201×5="
Replace-Text "156×7=" "366×2="
Replace-Text "972×4=" "548×8="
Replace-Text "136×7=" "866×5="

Replace-Text "673×7=" "527×2="
Replace-Text "536×5=" "435×9="
Replace-Text "497×6=" "261×3="
Replace-Text "874×9=" "607×8="
Replace-Text "674×5=" "876×4="

Replace-Text "669×9=" "112×7="
Replace-Text "461×8=" "476×7="
Replace-Text "567×3=" "333×5="
Replace-Text "787×5=" "838×7="
Replace-Text "131×7=" "132×2="

Replace-Text "887×5=" "868×6="
Replace-Text "756×5=" "430×3="
Replace-Text "276×8=" "590×8="
Replace-Text "956×2=" "610×9="
Replace-Text "520×3=" "908×3="
